# Commit: Calculating the relative error to the value from the first data
# source result.
#
# The "Relative" column (previously Difference / 2nd-Source) is recomputed
# as Difference / 1st-Source. Only the cells whose ratio actually changes
# (i.e. rows that have both a 1st- and 2nd-source value) are affected; the
# "added/removed row" cases where the ratio was already 1.0 stay the same.
#
# The "Relative" column is formatted as a 0.000% percentage, so after the
# numbers change, its best-fit column width also grows (to fit the new
# longest value, 486.667%) and converges to the same width on every sheet.

$wb = $excel.ActiveWorkbook

# sheetName -> column letter that holds the "Relative" value
$relativeColumn = @{
    "full"       = "G"
    "left"       = "G"
    "right"      = "G"
    "inner"      = "G"
    "diffs"      = "G"
    "append"     = "G"
    "multicols"  = "H"
}

# sheetName -> list of (cellRef, newValue) pairs to update
$updates = @{
    "full"      = @(@("G5", 0.9090909090909091), @("G8", 0.8333333333333334), @("G26", 4.866666666666666))
    "left"      = @(@("G5", 0.9090909090909091), @("G8", 0.8333333333333334), @("G20", 4.866666666666666))
    "right"     = @(@("G5", 0.9090909090909091), @("G8", 0.8333333333333334), @("G20", 4.866666666666666))
    "inner"     = @(@("G5", 0.9090909090909091), @("G8", 0.8333333333333334), @("G14", 4.866666666666666))
    "diffs"     = @(@("G5", 0.9090909090909091), @("G8", 0.8333333333333334), @("G21", 4.866666666666666))
    "append"    = @(@("G5", 0.9090909090909091), @("G8", 0.8333333333333334), @("G9", 4.866666666666666), @("G12", 0.9090909090909091), @("G15", 0.8333333333333334), @("G16", 4.866666666666666))
    "multicols" = @(@("H4", 0.9090909090909091), @("H7", 0.8333333333333334), @("H20", 4.866666666666666))
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($pair in $updates[$sheetName]) {
        $cellRef = $pair[0]
        $newValue = $pair[1]
        $ws.Range($cellRef).Value = $newValue
    }

    # Best-fit the "Relative" column for its new, longer formatted values.
    $colLetter = $relativeColumn[$sheetName]
    $ws.Range("$colLetter`1").EntireColumn.ColumnWidth = 12.8
}
